$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Shift the closing signature rows (22-23) down to (23-24) to make room ---
# Row 22/23 have no neighbours below used, so a plain row insert moves their
# content+styles down intact.
$ws.Rows.Item(22).Insert()

# --- 2) Duplicate row 17 (old second detail row, style "last row") into the
#        new row 18 slot, carrying its style + its current values along ---
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))

# --- 3) Duplicate row 16 (first detail row, style "middle row") into row 17 ---
# Row 16's current content (CC/1050958181/MERLYS GELIZ LLORENA/1904/32400/900000)
# becomes row 17's target content, while row 17 now also carries row16's style.
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))

# --- 4) Row 16 keeps its own (middle-row) style; update its values to the
#        second period for MERLYS GELIZ LLORENA (previously on old row 17) ---
$ws.Range("E16").Value = "1905"
$ws.Range("F16").Value = 36000

# --- 5) Fill in the brand-new worker row (row 18) ---
$ws.Range("C18").Value = "1007981090"
$ws.Range("D18").Value = "ALDAIR TURIZO RUIZ"
$ws.Range("E18").Value = "2503"
$ws.Range("F18").Value = 9490
$ws.Range("G18").Value = 1423500

# --- 6) Update the summary counters and total overdue amount ---
$ws.Range("E11").Value = 77890
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 3
